# Update gh-pages output data (漫展信息) across the 展览 / 演出 / 全部类型 sheets.
# Values below mirror a refreshed scrape: "want-to-go" counts (col F),
# a couple of min-ticket-price updates (col G), and one refreshed cover
# image URL (col I) on row 12 of 展览 / 全部类型.

$wb = $excel.ActiveWorkbook

# ---- 展览 ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 12
$ws.Range("F3").Value = 5249
$ws.Range("F4").Value = 12
$ws.Range("F5").Value = 7547
$ws.Range("F6").Value = 46
$ws.Range("F8").Value = 103
$ws.Range("F9").Value = 605
$ws.Range("F12").Value = 4358
$ws.Range("I12").Value = "//i1.hdslb.com/bfs/openplatform/202407/DSmxF8X51721900155471.jpeg"
$ws.Range("F13").Value = 1777
$ws.Range("F15").Value = 113
$ws.Range("F16").Value = 2944
$ws.Range("F20").Value = 524
$ws.Range("F21").Value = 456
$ws.Range("F22").Value = 467
$ws.Range("F23").Value = 323
$ws.Range("F24").Value = 109
$ws.Range("F25").Value = 1706
$ws.Range("F26").Value = 1208
$ws.Range("G26").Value = 50
$ws.Range("F28").Value = 1402
$ws.Range("F30").Value = 586
$ws.Range("F31").Value = 30
$ws.Range("F34").Value = 11
$ws.Range("F35").Value = 66
$ws.Range("F38").Value = 2993
$ws.Range("F39").Value = 711
$ws.Range("F40").Value = 33
$ws.Range("F41").Value = 105
$ws.Range("F43").Value = 58

# ---- 演出 ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 19

# ---- 全部类型 ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 12
$ws.Range("F3").Value = 5249
$ws.Range("F4").Value = 12
$ws.Range("F5").Value = 7547
$ws.Range("F6").Value = 46
$ws.Range("F8").Value = 103
$ws.Range("F9").Value = 605
$ws.Range("F12").Value = 4358
$ws.Range("I12").Value = "//i1.hdslb.com/bfs/openplatform/202407/DSmxF8X51721900155471.jpeg"
$ws.Range("F13").Value = 1777
$ws.Range("F15").Value = 113
$ws.Range("F16").Value = 2944
$ws.Range("F20").Value = 524
$ws.Range("F21").Value = 456
$ws.Range("F22").Value = 467
$ws.Range("F24").Value = 323
$ws.Range("F25").Value = 109
$ws.Range("F26").Value = 1706
$ws.Range("F27").Value = 1208
$ws.Range("G27").Value = 50
$ws.Range("F29").Value = 1402
$ws.Range("F31").Value = 586
$ws.Range("F32").Value = 30
$ws.Range("F33").Value = 517
$ws.Range("F35").Value = 11
$ws.Range("F36").Value = 66
$ws.Range("F39").Value = 2993
$ws.Range("F40").Value = 19
$ws.Range("F41").Value = 711
$ws.Range("F42").Value = 33
$ws.Range("F43").Value = 105
$ws.Range("F45").Value = 58
